$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 607.3333
$ws.Range("J2").Value = 714.6667
$ws.Range("L2").Value = 714.6667
$ws.Range("N2").Value = -940.6667
$ws.Range("H57").Value = 43499.5
$ws.Range("I57").Value = 43499.5
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 130498.5
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -129999.5
$ws.Range("N57").ClearContents()
$ws.Range("H70").Value = 10000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 10000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 30000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -30540
$ws.Range("H73").Value = 10000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 10000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 30000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -31872
$ws.Range("H137").Value = 2005.9
$ws.Range("I137").Value = 1121.2
$ws.Range("J137").Value = 2890.6
$ws.Range("K137").Value = 3363.6
$ws.Range("L137").Value = 8671.799999999999
$ws.Range("M137").Value = -813.6000000000004
$ws.Range("N137").Value = -13771.8
$ws.Range("H139").Value = 130000
$ws.Range("J139").Value = 130000
$ws.Range("L139").Value = 130000
$ws.Range("N139").Value = -140280
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 699
$ws.Range("I4").Value = 699
$ws.Range("K4").Value = 699
$ws.Range("M4").Value = -583
$ws.Range("H5").Value = 181
$ws.Range("J5").Value = 255
$ws.Range("L5").Value = 255
$ws.Range("N5").Value = -479
$ws.Range("H6").Value = 18891000
$ws.Range("I6").Value = 18415714
$ws.Range("J6").Value = 20000000
$ws.Range("K6").Value = 18415714
$ws.Range("L6").Value = 20000000
$ws.Range("M6").Value = -18415541
$ws.Range("N6").Value = -20000346
$ws.Range("H97").Value = 1133.0769
$ws.Range("I97").Value = 700.8
$ws.Range("K97").Value = 700.8
$ws.Range("M97").Value = -204.8
$ws.Range("H132").Value = 1748.186
$ws.Range("I132").Value = 1776.2051
$ws.Range("K132").Value = 5328.615299999999
$ws.Range("M132").Value = -2798.615299999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 181
$ws.Range("J4").Value = 255
$ws.Range("L4").Value = 255
$ws.Range("N4").Value = -485
$ws.Range("H60").Value = 92750
$ws.Range("J60").Value = 92750
$ws.Range("L60").Value = 92750
$ws.Range("N60").Value = -93948
$ws.Range("H74").Value = 144666.67
$ws.Range("J74").Value = 144666.67
$ws.Range("L74").Value = 144666.67
$ws.Range("N74").Value = -146538.67
$ws.Range("H77").Value = 144666.67
$ws.Range("J77").Value = 144666.67
$ws.Range("L77").Value = 434000.01
$ws.Range("N77").Value = -443360.01
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1667728.6
$ws.Range("I3").Value = 3333791.8
$ws.Range("J3").Value = 1665.6666
$ws.Range("K3").Value = 3333791.8
$ws.Range("L3").Value = 1665.6666
$ws.Range("M3").Value = -3333678.8
$ws.Range("N3").Value = -1891.6666
$ws.Range("H7").Value = 74.46666999999999
$ws.Range("I7").Value = 51.916668
$ws.Range("K7").Value = 51.916668
$ws.Range("M7").Value = 61.083332
$ws.Range("H108").Value = 55513.25
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 55513.25
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 55513.25
$ws.Range("M108").ClearContents()
$ws.Range("N108").Value = -63193.25
$ws.Range("H109").Value = 47380.832
$ws.Range("J109").Value = 47380.832
$ws.Range("L109").Value = 47380.832
$ws.Range("N109").Value = -49460.832
$ws.Range("H132").Value = 4048.7646
$ws.Range("I132").Value = 3801.4167
$ws.Range("J132").Value = 4642.4
$ws.Range("K132").Value = 11404.2501
$ws.Range("L132").Value = 13927.2
$ws.Range("M132").Value = -8874.250100000001
$ws.Range("N132").Value = -18987.2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 3232.6667
$ws.Range("J22").Value = 3399.5
$ws.Range("L22").Value = 10198.5
$ws.Range("N22").Value = -10536.5
$ws.Range("H27").Value = 3232.6667
$ws.Range("J27").Value = 3399.5
$ws.Range("L27").Value = 10198.5
$ws.Range("N27").Value = -10402.5
$ws.Range("H29").Value = 137.14285
$ws.Range("I29").Value = 13.333333
$ws.Range("J29").Value = 230
$ws.Range("K29").Value = 39.999999
$ws.Range("L29").Value = 690
$ws.Range("M29").Value = 237.000001
$ws.Range("N29").Value = -1244
$ws.Range("H122").Value = 870.3333
$ws.Range("I122").Value = 833.7143
$ws.Range("K122").Value = 7503.428699999999
$ws.Range("M122").Value = -5053.428699999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3319.6667
$ws.Range("I126").Value = 3061.1538
$ws.Range("K126").Value = 9183.4614
$ws.Range("M126").Value = -6713.4614
$ws.Range("H132").Value = 46982
$ws.Range("I132").Value = 73781.42999999999
$ws.Range("J132").Value = 5294
$ws.Range("K132").Value = 221344.29
$ws.Range("L132").Value = 15882
$ws.Range("M132").Value = -218814.29
$ws.Range("N132").Value = -20942
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 863.6
$ws.Range("I22").Value = 841
$ws.Range("J22").Value = 897.5
$ws.Range("K22").Value = 841
$ws.Range("L22").Value = 897.5
$ws.Range("M22").Value = -546
$ws.Range("N22").Value = -1487.5
$ws.Range("H27").Value = 863.6
$ws.Range("I27").Value = 841
$ws.Range("J27").Value = 897.5
$ws.Range("K27").Value = 841
$ws.Range("L27").Value = 897.5
$ws.Range("M27").Value = -734
$ws.Range("N27").Value = -1111.5
$ws.Range("H40").Value = 6548.8
$ws.Range("I40").Value = 6187.25
$ws.Range("K40").Value = 6187.25
$ws.Range("M40").Value = -6051.25
$ws.Range("H122").Value = 4241
$ws.Range("J122").Value = 4474.5
$ws.Range("L122").Value = 13423.5
$ws.Range("N122").Value = -18323.5
$ws.Range("H140").Value = 113331.664
$ws.Range("J140").Value = 69997.5
$ws.Range("L140").Value = 69997.5
$ws.Range("N140").Value = -80357.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 22334.666
$ws.Range("J3").Value = 57004
$ws.Range("L3").Value = 57004
$ws.Range("N3").Value = -57232
$ws.Range("H14").Value = 1404.8334
$ws.Range("I14").Value = 2701.3333
$ws.Range("J14").Value = 108.333336
$ws.Range("K14").Value = 2701.3333
$ws.Range("L14").Value = 108.333336
$ws.Range("M14").Value = -2533.3333
$ws.Range("N14").Value = -444.333336
$ws.Range("H126").Value = 3614.739
$ws.Range("I126").Value = 1976.6
$ws.Range("K126").Value = 5929.799999999999
$ws.Range("M126").Value = -3459.799999999999
$ws.Range("H132").Value = 2405.5
$ws.Range("I132").Value = 2405.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7216.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4686.5
$ws.Range("N132").ClearContents()
